$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "MSG: None`n`nMSG: The decision to select a movie for Friday has been recorded as `"no decision.`"`n"
$ws.Range("D2").Value = "no_decision, "
$ws.Range("C3").Value = "MSG: None`n`nMSG: The decision has been recorded, and the movie `"Barbie`" will be acquired for showing on Friday.`n"
$ws.Range("C4").Value = "MSG: None`n`nMSG: The decision to acquire the rights for both movies has been recorded.`n"
$ws.Range("D4").Value = "both_movies, "
$ws.Range("C5").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for `"Barbie.`"`n"
$ws.Range("C6").Value = "MSG: None`n`nMSG: The decision about what movie to show on Friday was not reached.`n"
$ws.Range("D6").Value = "no_decision, "
$ws.Range("C7").Value = "MSG: None`n`nMSG: The decision for a movie to be shown on Friday could not be reached, resulting in no decision being made.`n"
$ws.Range("D7").Value = "no_decision, "
$ws.Range("C8").Value = "MSG: None`n`nMSG: I have recorded the decision to acquire the rights for the movie `"Barbie.`"`n"
$ws.Range("C9").Value = "MSG: None`n`nMSG: The conversation ended without a definitive decision on what movie to show on Friday, so I have called the no_decision function.`n"
$ws.Range("D9").Value = "no_decision, "
$ws.Range("C10").Value = "MSG: None`n`nMSG: The decision has been made that there is no clear choice for the movie to be shown on Friday.`n"
$ws.Range("D10").Value = "no_decision, "
$ws.Range("C11").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for `"Barbie.`"`n"
$ws.Range("C12").Value = "MSG: None`n`nMSG: The decision-making process did not result in an agreement on which movie to show on Friday, so the conclusion is that there is no decision.`n"
$ws.Range("D12").Value = "no_decision, "
$ws.Range("C13").Value = "MSG: None`n`nMSG: The decision regarding the movie for Friday was not made, and the conversation ended without a definitive plan.`n"
$ws.Range("D13").Value = "no_decision, "
$ws.Range("C14").Value = "MSG: None`n`nMSG: The decision has been recorded to acquire the rights for `"Barbie.`"`n"
$ws.Range("C15").Value = "MSG: None`n`nMSG: The decision regarding which movie to show on Friday could not be made, so the choice remains undecided.`n"
$ws.Range("D15").Value = "no_decision, "
$ws.Range("C16").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for the movie `"Barbie.`"`n"
$ws.Range("C17").Value = "MSG: None`n`nMSG: The decision about which movie to show on Friday resulted in no agreement.`n"
$ws.Range("D17").Value = "no_decision, "
$ws.Range("C18").Value = "MSG: None`n`nMSG: Based on the given criteria, there was no explicit decision made about the movie to be shown on Friday, so the appropriate action is to call the no_decision function.`n"
$ws.Range("D18").Value = "no_decision, "
$ws.Range("C19").Value = "MSG: None`n`nMSG: The decision to acquire `"Barbie`" rights has been made successfully.`n"
$ws.Range("C20").Value = "MSG: None`n`nMSG: The decision has been recorded, and no movie was selected for Friday.`n"
$ws.Range("D20").Value = "no_decision, "
$ws.Range("C21").Value = "MSG: None`n`nMSG: The decision has been recorded, and no movie was definitively chosen for this Friday.`n"
$ws.Range("D21").Value = "no_decision, "
$ws.Range("C22").Value = "MSG: None`n`nMSG: The decision has been recorded as no movie being selected for Friday.`n"
$ws.Range("C23").Value = "MSG: None`n`nMSG: The decision has been recorded as no choice made regarding Friday's movie.`n"
$ws.Range("D23").Value = "no_decision, "
$ws.Range("C24").Value = "MSG: None`n`nMSG: The decision process resulted in no clear agreement regarding the movie to show on Friday.`n"
$ws.Range("D24").Value = "no_decision, "
$ws.Range("C25").Value = "MSG: None`n`nMSG: The decision to acquire the rights for both movies has been recorded successfully.`n"
$ws.Range("D25").Value = "both_movies, "
$ws.Range("C26").Value = "MSG: None`n`nMSG: The decision has been recorded as no decision regarding the movie for Friday.`n"
$ws.Range("D26").Value = "no_decision, "
$ws.Range("C27").Value = "MSG: None`n`nMSG: The decision has been recorded, indicating that no movie has been selected for Friday.`n"
$ws.Range("D27").Value = "no_decision, "
$ws.Range("C28").Value = "MSG: None`n`nMSG: The decision to acquire the rights for both movies has been recorded.`n"
$ws.Range("D28").Value = "both_movies, "
$ws.Range("C29").Value = "MSG: None`n`nMSG: The decision process concluded without a definitive choice for the movie on Friday, so no action was taken to acquire rights for any film.`n"
$ws.Range("D29").Value = "no_decision, "
$ws.Range("C30").Value = "MSG: None`n`nMSG: The decision has been recorded, and the movie `"Barbie`" will be acquired for screening.`n"
$ws.Range("C31").Value = "MSG: None`n`nMSG: The decision about what movie to show on Friday has been recorded as no decision.`n"
$ws.Range("D31").Value = "no_decision, "
$ws.Range("C32").Value = "MSG: None`n`nMSG: The decision has been recorded as no decision regarding the movie for Friday.`n"
$ws.Range("D32").Value = "no_decision, "
$ws.Range("C33").Value = "MSG: None`n`nMSG: The decision has been recorded, and no definitive choice for the movie to show on Friday was made.`n"
$ws.Range("D33").Value = "no_decision, "
$ws.Range("C34").Value = "MSG: None`n`nMSG: The decision about what movie to show on Friday has not been made, indicating no selection was reached.`n"
$ws.Range("D34").Value = "no_decision, "
$ws.Range("C35").Value = "MSG: None`n`nMSG: The decision regarding the movie for Friday resulted in no agreement, as confirmed by the call to the no_decision function.`n"
$ws.Range("D35").Value = "no_decision, "
$ws.Range("C36").Value = "MSG: None`n`nMSG: The conversation ended without a decision about which movie to show on Friday.`n"
$ws.Range("D36").Value = "no_decision, "
$ws.Range("C37").Value = "MSG: None`n`nMSG: The decision has been recorded as `"no decision.`"`n"
$ws.Range("D37").Value = "no_decision, "
$ws.Range("C38").Value = "MSG: None`n`nMSG: I have recorded the decision to acquire the rights for `"Barbie`" for the movie to be shown on Friday.`n"
$ws.Range("C39").Value = "MSG: None`n`nMSG: No decision was made regarding the movie to be shown on Friday.`n"
$ws.Range("D39").Value = "no_decision, "
$ws.Range("C40").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for the movie `"Barbie.`"`n"
$ws.Range("C41").Value = "MSG: None`n`nMSG: The decision has been recorded as no choice could be made regarding Friday's movie.`n"
$ws.Range("D41").Value = "no_decision, "
$ws.Range("C42").Value = "MSG: None`n`nMSG: The decision has been recorded, and the movie `"Barbie`" will be acquired for screening on Friday.`n"
$ws.Range("C43").Value = "MSG: None`n`nMSG: The decision has been recorded with no movie selected for Friday.`n"
$ws.Range("D43").Value = "no_decision, "
$ws.Range("C44").Value = "MSG: None`n`nMSG: The rights to `"Barbie`" have been successfully acquired.`n"
$ws.Range("C45").Value = "MSG: None`n`nMSG: The decision to show `"Barbie`" has been made.`n"
$ws.Range("C46").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for `"Barbie.`"`n"
$ws.Range("C47").Value = "MSG: None`n`nMSG: The decision about which movie to show on Friday was not reached, as the committee did not come to a conclusive choice.`n"
$ws.Range("D47").Value = "no_decision, "
$ws.Range("C48").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for `"Barbie.`"`n"
$ws.Range("C49").Value = "MSG: None`n`nMSG: No decision was made regarding the movie to play on Friday.`n"
$ws.Range("C50").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for `"Oppenheimer.`"`n"
$ws.Range("C51").Value = "MSG: None`n`nMSG: The decision about what movie to show on Friday could not be determined.`n"
$ws.Range("D51").Value = "no_decision, "
$ws.Range("C52").Value = "MSG: None`n`nMSG: The decision about what movie to show on Friday was not made.`n"
$ws.Range("D52").Value = "no_decision, "
$ws.Range("C53").Value = "MSG: None`n`nMSG: The result indicates that there was no decision made regarding the movie to show on Friday.`n"
$ws.Range("D53").Value = "no_decision, "
$ws.Range("C54").Value = "MSG: None`n`nMSG: The decision is to record no decision regarding which movie to play on Friday.`n"
$ws.Range("D54").Value = "no_decision, "
$ws.Range("C55").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for `"Barbie`" for Friday's showing.`n"
$ws.Range("C56").Value = "MSG: None`n`nMSG: The decision to acquire the rights for `"Barbie`" has been successfully recorded.`n"
$ws.Range("C57").Value = "MSG: None`n`nMSG: The rights to both movies have been successfully acquired.`n"
$ws.Range("D57").Value = "both_movies, "
$ws.Range("C58").Value = "MSG: None`n`nMSG: The decision regarding the movie selection has concluded with no movie chosen. If you need any further assistance or information, feel free to ask!`n"
$ws.Range("D58").Value = "no_decision, "
$ws.Range("C59").Value = "MSG: None`n`nMSG: The decision regarding Friday's movie has concluded without a selection.`n"
$ws.Range("D59").Value = "no_decision, "
$ws.Range("C60").Value = "MSG: None`n`nMSG: The decision has been recorded, and no movie will be acquired at this time.`n"
$ws.Range("D60").Value = "no_decision, "
$ws.Range("C61").Value = "MSG: None`n`nMSG: The decision has been recorded as no decision regarding the movie for Friday was made.`n"
$ws.Range("D61").Value = "no_decision, "
$ws.Range("C62").Value = "MSG: None`n`nMSG: The rights to both movies have been successfully acquired.`n"
$ws.Range("D62").Value = "both_movies, "
$ws.Range("C63").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for `"Barbie.`"`n"
$ws.Range("C64").Value = "MSG: None`n`nMSG: The decision has been recorded as no movie was selected for Friday.`n"
$ws.Range("D64").Value = "no_decision, "
$ws.Range("C65").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for `"Barbie`" to be shown on Friday.`n"
$ws.Range("C66").Value = "MSG: None`n`nMSG: The conversation did not result in a clearly defined decision regarding which movie will be shown on Friday.`n"
$ws.Range("D66").Value = "no_decision, "
$ws.Range("C67").Value = "MSG: None`n`nMSG: The decision has been made to show `"Barbie`" on Friday.`n"
$ws.Range("C68").Value = "MSG: None`n`nMSG: The decision has been recorded as no decision being made regarding the movie to show on Friday.`n"
$ws.Range("D68").Value = "no_decision, "
